$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current row 361 ("Packham's Triumph" / Primera,
# date 44432). This pushes the existing rows 361-362 down to 365-366 untouched, and
# leaves 4 fresh rows (361-364) to populate with the new weekly entries.
$ws.Range("A361:A364").EntireRow.Insert()

# Common fixed fields shared by every row in this block.
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$prodId    = 100104
$producto  = "Frutos de pepita"
$catId     = 100104005
$categoria = "Pera"
$unidad    = "$/caja 16 kilos empedrada"
$origen    = "Región de O'Higgins"
$kgUnidad  = 16

function Set-PeraRow($row, $fecha, $variedad, $calidad, $volumen, $pmin, $pmax, $pprom, $pkg) {
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $prodId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $pmin
    $ws.Cells.Item($row, 15).Value = $pmax
    $ws.Cells.Item($row, 16).Value = $pprom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $pkg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

# New week (44656) entries.
Set-PeraRow 361 44656 "Abate Fettel"        "Primera" 100 8000 9000  8500 531
Set-PeraRow 362 44656 "Abate Fettel"        "Segunda" 50  7000 7000  7000 438
Set-PeraRow 363 44656 "Packham's Triumph"   "Primera" 150 9000 11000 9667 604
Set-PeraRow 364 44656 "Packham's Triumph"   "Segunda" 100 8000 8000  8000 500
